$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column BQ: header label "11-sep" (continues the date-label header row),
# plus one numeric value per data row (2-18), mirroring the existing BP column.
$ws.Range("BQ1").Value = "11-sep"

$ws.Range("BQ2").Value  = 0
$ws.Range("BQ3").Value  = 18.871834183711847
$ws.Range("BQ4").Value  = 14.668806985304631
$ws.Range("BQ5").Value  = 16.683398571673791
$ws.Range("BQ6").Value  = 0
$ws.Range("BQ7").Value  = 13.566582448166319
$ws.Range("BQ8").Value  = 11.688227970745805
$ws.Range("BQ9").Value  = 15.710353396536831
$ws.Range("BQ10").Value = 15.26541066913107
$ws.Range("BQ11").Value = 11.056731268974424
$ws.Range("BQ12").Value = 0
$ws.Range("BQ13").Value = 11.794881493834833
$ws.Range("BQ14").Value = 0
$ws.Range("BQ15").Value = 0
$ws.Range("BQ16").Value = 10.264797993271037
$ws.Range("BQ17").Value = 0
$ws.Range("BQ18").Value = 0

# Cursor/selection moved to BS4 in the saved file.
$ws.Range("BS4").Select()
